$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(19, 9).Value = 'b'
$ws.Cells.Item(19, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(30, 9).Value = '%'
$ws.Cells.Item(30, 10).Value = 'Uninterpretable'
$ws.Cells.Item(45, 9).Value = 'sd'
$ws.Cells.Item(45, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(55, 9).Value = 'sd'
$ws.Cells.Item(55, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(76, 9).Value = 'sd'
$ws.Cells.Item(76, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(78, 9).Value = 'sd'
$ws.Cells.Item(78, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(92, 9).Value = 'aa'
$ws.Cells.Item(92, 10).Value = 'Agree/Accept'
$ws.Cells.Item(103, 9).Value = 'sd'
$ws.Cells.Item(103, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(125, 9).Value = 'sd'
$ws.Cells.Item(125, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(128, 9).Value = 'sv'
$ws.Cells.Item(128, 10).Value = 'Statement-opinion'
$ws.Cells.Item(136, 9).Value = 'sv'
$ws.Cells.Item(136, 10).Value = 'Statement-opinion'
$ws.Cells.Item(144, 9).Value = 'sv'
$ws.Cells.Item(144, 10).Value = 'Statement-opinion'
$ws.Cells.Item(163, 9).Value = 'sv'
$ws.Cells.Item(163, 10).Value = 'Statement-opinion'
$ws.Cells.Item(164, 9).Value = 'sv'
$ws.Cells.Item(164, 10).Value = 'Statement-opinion'
$ws.Cells.Item(165, 9).Value = 'sd'
$ws.Cells.Item(165, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(170, 9).Value = 'sv'
$ws.Cells.Item(170, 10).Value = 'Statement-opinion'
$ws.Cells.Item(175, 9).Value = 'ba'
$ws.Cells.Item(175, 10).Value = 'Appreciation'
$ws.Cells.Item(178, 9).Value = 'sd'
$ws.Cells.Item(178, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(181, 9).Value = 'sd'
$ws.Cells.Item(181, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(185, 9).Value = '%'
$ws.Cells.Item(185, 10).Value = 'Uninterpretable'
$ws.Cells.Item(186, 9).Value = 'sv'
$ws.Cells.Item(186, 10).Value = 'Statement-opinion'
$ws.Cells.Item(201, 9).Value = 'sd'
$ws.Cells.Item(201, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(209, 9).Value = 'sd'
$ws.Cells.Item(209, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(235, 9).Value = 'ba'
$ws.Cells.Item(235, 10).Value = 'Appreciation'
$ws.Cells.Item(238, 9).Value = 'sd'
$ws.Cells.Item(238, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(240, 9).Value = 'aa'
$ws.Cells.Item(240, 10).Value = 'Agree/Accept'
$ws.Cells.Item(243, 9).Value = 'sv'
$ws.Cells.Item(243, 10).Value = 'Statement-opinion'
$ws.Cells.Item(248, 9).Value = 'sd'
$ws.Cells.Item(248, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(274, 9).Value = 'sv'
$ws.Cells.Item(274, 10).Value = 'Statement-opinion'
$ws.Cells.Item(279, 9).Value = 'sv'
$ws.Cells.Item(279, 10).Value = 'Statement-opinion'
$ws.Cells.Item(286, 9).Value = '%'
$ws.Cells.Item(286, 10).Value = 'Uninterpretable'
$ws.Cells.Item(295, 9).Value = 'sd'
$ws.Cells.Item(295, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(309, 9).Value = 'b'
$ws.Cells.Item(309, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(312, 9).Value = 'sd'
$ws.Cells.Item(312, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(316, 9).Value = 'sd'
$ws.Cells.Item(316, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(342, 9).Value = 'sd'
$ws.Cells.Item(342, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(350, 9).Value = 'sd'
$ws.Cells.Item(350, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(353, 9).Value = 'sd'
$ws.Cells.Item(353, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(354, 9).Value = '%'
$ws.Cells.Item(354, 10).Value = 'Uninterpretable'
$ws.Cells.Item(356, 9).Value = '%'
$ws.Cells.Item(356, 10).Value = 'Uninterpretable'
$ws.Cells.Item(366, 9).Value = 'ba'
$ws.Cells.Item(366, 10).Value = 'Appreciation'
$ws.Cells.Item(369, 9).Value = 'ba'
$ws.Cells.Item(369, 10).Value = 'Appreciation'
$ws.Cells.Item(389, 9).Value = 'b'
$ws.Cells.Item(389, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(399, 9).Value = '%'
$ws.Cells.Item(399, 10).Value = 'Uninterpretable'
$ws.Cells.Item(414, 9).Value = 'sd'
$ws.Cells.Item(414, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(425, 9).Value = 'aa'
$ws.Cells.Item(425, 10).Value = 'Agree/Accept'
$ws.Cells.Item(427, 9).Value = 'b'
$ws.Cells.Item(427, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(433, 9).Value = 'sv'
$ws.Cells.Item(433, 10).Value = 'Statement-opinion'
$ws.Cells.Item(442, 9).Value = 'sd'
$ws.Cells.Item(442, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(443, 9).Value = 'sd'
$ws.Cells.Item(443, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(444, 9).Value = 'sv'
$ws.Cells.Item(444, 10).Value = 'Statement-opinion'
$ws.Cells.Item(462, 9).Value = 'sd'
$ws.Cells.Item(462, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(472, 9).Value = 'sd'
$ws.Cells.Item(472, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(488, 9).Value = 'ba'
$ws.Cells.Item(488, 10).Value = 'Appreciation'
$ws.Cells.Item(493, 9).Value = '%'
$ws.Cells.Item(493, 10).Value = 'Uninterpretable'
$ws.Cells.Item(501, 9).Value = '%'
$ws.Cells.Item(501, 10).Value = 'Uninterpretable'
$ws.Cells.Item(517, 9).Value = 'sv'
$ws.Cells.Item(517, 10).Value = 'Statement-opinion'
$ws.Cells.Item(537, 9).Value = 'sv'
$ws.Cells.Item(537, 10).Value = 'Statement-opinion'
$ws.Cells.Item(540, 9).Value = 'sd'
$ws.Cells.Item(540, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(543, 9).Value = 'ba'
$ws.Cells.Item(543, 10).Value = 'Appreciation'
$ws.Cells.Item(548, 9).Value = '%'
$ws.Cells.Item(548, 10).Value = 'Uninterpretable'
$ws.Cells.Item(566, 9).Value = '%'
$ws.Cells.Item(566, 10).Value = 'Uninterpretable'
$ws.Cells.Item(574, 9).Value = 'sv'
$ws.Cells.Item(574, 10).Value = 'Statement-opinion'
$ws.Cells.Item(575, 9).Value = 'sd'
$ws.Cells.Item(575, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(580, 9).Value = 'b'
$ws.Cells.Item(580, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(582, 9).Value = 'sv'
$ws.Cells.Item(582, 10).Value = 'Statement-opinion'
$ws.Cells.Item(593, 9).Value = 'sv'
$ws.Cells.Item(593, 10).Value = 'Statement-opinion'
$ws.Cells.Item(601, 9).Value = 'sv'
$ws.Cells.Item(601, 10).Value = 'Statement-opinion'
$ws.Cells.Item(605, 9).Value = 'sv'
$ws.Cells.Item(605, 10).Value = 'Statement-opinion'
$ws.Cells.Item(625, 9).Value = 'sd'
$ws.Cells.Item(625, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(631, 9).Value = '%'
$ws.Cells.Item(631, 10).Value = 'Uninterpretable'
$ws.Cells.Item(645, 9).Value = 'sv'
$ws.Cells.Item(645, 10).Value = 'Statement-opinion'
$ws.Cells.Item(647, 9).Value = 'sv'
$ws.Cells.Item(647, 10).Value = 'Statement-opinion'
$ws.Cells.Item(662, 9).Value = 'sv'
$ws.Cells.Item(662, 10).Value = 'Statement-opinion'
$ws.Cells.Item(672, 9).Value = 'sv'
$ws.Cells.Item(672, 10).Value = 'Statement-opinion'
$ws.Cells.Item(674, 9).Value = 'sv'
$ws.Cells.Item(674, 10).Value = 'Statement-opinion'
$ws.Cells.Item(687, 9).Value = '%'
$ws.Cells.Item(687, 10).Value = 'Uninterpretable'
$ws.Cells.Item(692, 9).Value = 'sd'
$ws.Cells.Item(692, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(694, 9).Value = 'sv'
$ws.Cells.Item(694, 10).Value = 'Statement-opinion'
$ws.Cells.Item(696, 9).Value = 'aa'
$ws.Cells.Item(696, 10).Value = 'Agree/Accept'
$ws.Cells.Item(697, 9).Value = 'aa'
$ws.Cells.Item(697, 10).Value = 'Agree/Accept'
$ws.Cells.Item(699, 9).Value = 'sv'
$ws.Cells.Item(699, 10).Value = 'Statement-opinion'
$ws.Cells.Item(700, 9).Value = 'sv'
$ws.Cells.Item(700, 10).Value = 'Statement-opinion'
$ws.Cells.Item(701, 9).Value = 'sd'
$ws.Cells.Item(701, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(709, 9).Value = 'aa'
$ws.Cells.Item(709, 10).Value = 'Agree/Accept'
